$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Toyota Camry VIN (choice)
$ws.Range("A2").Value = "4T1BE30K&6"
$ws.Range("B2").Value = "SYMBOL_2000_CHOICE"
$ws.Range("C2").Value = 2006
$ws.Range("D2").Value = "TOYT"
$ws.Range("E2").Value = "TOYOTA"
$ws.Range("F2").Value = "CAMRY"
$ws.Range("G2").Value = "CAMRY LE/XLE/SE"
$ws.Range("H2").Value = 20000
$ws.Range("I2").Value = "4D SED"
$ws.Range("J2").Value = "SEDAN 4 DOOR"
$ws.Range("M2").Value = "SED"
$ws.Range("N2").Value = "2.4L L4"
$ws.Range("O2").Value = 4
$ws.Range("R2").Value = "2WD"
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = "000S"
$ws.Range("U2").Value = "FRONT, HEAD & SIDE AIRBAGS"
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = "4 WHEEL STANDARD"
$ws.Range("X2").Value = "STD"
$ws.Range("Y2").Value = "H-IMMOBILIZER/ALARM"
$ws.Range("Z2").Value = "I"
$ws.Range("AA2").Value = 11
$ws.Range("AB2").Value = 11
$ws.Range("AC2").Value = "A"
$ws.Range("AD2").Value = "N"

# Row 3 - Toyota Camry VIN (test choice tier)
$ws.Range("A3").Value = "4T1BE30K&6"
$ws.Range("B3").Value = "SYMBOL_2000_CHOICE_T"
$ws.Range("C3").Value = 2006
$ws.Range("D3").Value = "TEST"
$ws.Range("E3").Value = "TEST"
$ws.Range("F3").Value = "TEST"
$ws.Range("G3").Value = "TEST"
$ws.Range("H3").Value = 20000
$ws.Range("I3").Value = "TEST"
$ws.Range("J3").Value = "TEST"
$ws.Range("K3").Value = "TEST"
$ws.Range("L3").Value = "TEST"
$ws.Range("M3").Value = "SED"
$ws.Range("N3").Value = "2.4L L4"
$ws.Range("O3").Value = 4
$ws.Range("R3").Value = "2WD"
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = "000S"
$ws.Range("U3").Value = "FRONT, HEAD & SIDE AIRBAGS"
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = "4 WHEEL STANDARD"
$ws.Range("X3").Value = "STD"
$ws.Range("Y3").Value = "H-IMMOBILIZER/ALARM"
$ws.Range("Z3").Value = "I"
$ws.Range("AA3").Value = 12
$ws.Range("AB3").Value = 13
$ws.Range("AC3").Value = "A"
$ws.Range("AD3").Value = "N"

# Update the visible window / selection to match the saved view state
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 20
$ws.Range("AA15").Select()
